$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Invoice # changes from 6 to 8
$ws.Range("F3").Value = "8"

# Customer address correction
$ws.Range("A11").Value = "102, Whitestone Veroso, Bangalore 560049"

# Phone number correction (customer contact / phone)
$ws.Range("A13").Value = "9900019361"
$ws.Range("D16").Value = "9900019361"

# Line item 1 (row 19) updates
$ws.Range("B19").Value = "Test Book73"
$ws.Range("C19").Value = "2"
$ws.Range("D19").Value = "10"
$ws.Range("F19").Value = "20"

# Line item 2 (row 20) - previously blank, now populated
$ws.Range("A20").Value = "2"
$ws.Range("B20").Value = "Test Book78"
$ws.Range("C20").Value = "2"
$ws.Range("D20").Value = "10"
$ws.Range("E20").Value = "0.0"
$ws.Range("F20").Value = "20"

# Line item 3 (row 21) - previously blank, now populated
$ws.Range("A21").Value = "3"
$ws.Range("B21").Value = "Test Book76"
$ws.Range("C21").Value = "3"
$ws.Range("D21").Value = "50"
$ws.Range("E21").Value = "0.0"
$ws.Range("F21").Value = "150"

# Updated total
$ws.Range("F29").Value = "190"
